# Auto-generated script to update cryptos.xlsx per the target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.189.26"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "2.912.69"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'364.38"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").Value = "'103.51"
$ws.Range("E6").Value = "  -5.64%  "
$ws.Range("E7").Value = "  -4.57%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -6.98%  "
$ws.Range("D10").Value = "'36.96"
$ws.Range("E10").Value = "  -5.09%  "
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("E12").Value = "  -4.09%  "
$ws.Range("D13").Value = "'18.44"
$ws.Range("E13").Value = "  -5.43%  "
$ws.Range("D14").Value = "3.372.85"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "'7.34"
$ws.Range("E15").Value = "  -5.60%  "
$ws.Range("D16").Value = "2.918.70"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "'0.954"
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").Value = "51.142.23"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("E19").Value = "  -4.03%  "
$ws.Range("E20").Value = "  -3.90%  "
$ws.Range("D21").Value = "'13.00"
$ws.Range("E21").Value = "  -6.57%  "
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").Value = "'68.19"
$ws.Range("E23").Value = "  -3.35%  "
$ws.Range("D24").Value = "'259.86"
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("D25").Value = "'2.69"
$ws.Range("E25").Value = "  -4.61%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'4.32"
$ws.Range("E26").Value = "  +3.94%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.173"
$ws.Range("E27").Value = "  -5.97%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'25.87"
$ws.Range("E29").Value = "  -3.91%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'7.18"
$ws.Range("E30").Value = "  -6.08%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.105"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "'6.16"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'9.95"
$ws.Range("E33").Value = "  -5.14%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "'2.14"
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'35.07"
$ws.Range("E35").Value = "  -6.70%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'50.72"
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0423"
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'2.80"
$ws.Range("E39").Value = "  +3.04%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'3.14"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "'16.97"
$ws.Range("E41").Value = "  -6.98%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.86"
$ws.Range("E42").Value = "  -6.67%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.114"
$ws.Range("E43").Value = "  -4.93%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'22.26"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'119.12"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'2.11"
$ws.Range("E46").Value = "  -2.90%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.062.12"
$ws.Range("E47").Value = "  -3.13%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'3.19"
$ws.Range("E48").Value = "  -7.95%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'2.26"
$ws.Range("E49").Value = "  -8.40%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "3.198.18"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "'0.237"
$ws.Range("E51").Value = "  -5.81%  "

Write-Host "Applied 143 cell updates to cryptos sheet"
